# The "제도유형" (E) column previously held the repeated label
# "소단위전공과정, 마이크로디그리" for every microdegree row. The upload
# simplifies that label to just "마이크로디그리" for all data rows (2-38).
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("E2:E38").Value = "마이크로디그리"

# Restore the cursor/selection to cell B3, matching the saved workbook state.
$ws.Range("B3").Select()
